# "add s-shaped curves for the future and add comet 1"
#
# Rescale the ICAO target (col C) and IATA target (col D) projection
# curves for years 2005-2050 (rows 39-84) by a constant growth factor.
# Both series are simple geometric-decay projections (each year = prior
# year * 0.98 for col C, * 0.985 for col D); the edit shifts the whole
# curve upward by multiplying every existing value by the same constant
# factor, preserving the underlying decay shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$factor = 1.0075153088370112

for ($r = 39; $r -le 84; $r++) {
    $oldC = $ws.Range("C$r").Value()
    $oldD = $ws.Range("D$r").Value()
    $ws.Range("C$r").Value = $oldC * $factor
    $ws.Range("D$r").Value = $oldD * $factor
}
